# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Swap the Menorca / Ceuta rows (city name + stats), so that Ceuta now
# appears before Menorca in the shared-string table / data listing.
$ws.Range("A59").Value = "Ceuta"
$ws.Range("B59").Value = 16
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 16

$ws.Range("A60").Value = "Menorca"
$ws.Range("B60").Value = 15
$ws.Range("C60").Value = 18
$ws.Range("D60").Value = 13

# Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 10:42"
